$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.749.07'
$ws.Range("E2").Value = '  -6.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.224.40'
$ws.Range("E3").Value = '  -7.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.52'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.85'
$ws.Range("E6").Value = '  -13.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.573'
$ws.Range("E7").Value = '  -9.80%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  -10.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.13'
$ws.Range("E10").Value = '  -11.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.69'
$ws.Range("E11").Value = '  -4.52%  '
$ws.Range("E12").Value = '  -10.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.64'
$ws.Range("E13").Value = '  -12.52%  '
$ws.Range("E14").Value = '  -4.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.874'
$ws.Range("E15").Value = '  -13.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.570.60'
$ws.Range("E16").Value = '  -6.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.96'
$ws.Range("E17").Value = '  -11.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.215.28'
$ws.Range("E18").Value = '  -7.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.634.11'
$ws.Range("E19").Value = '  -6.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.01'
$ws.Range("E20").Value = '  +4.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.65'
$ws.Range("E21").Value = '  -11.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0949'
$ws.Range("E22").Value = '  -12.56%  '
$ws.Range("E23").Value = '  -7.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.85'
$ws.Range("E24").Value = '  -13.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '235.12'
$ws.Range("E25").Value = '  -11.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.12'
$ws.Range("E26").Value = '  -10.09%  '
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.20'
$ws.Range("E28").Value = '  -9.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.56'
$ws.Range("E29").Value = '  -13.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.18'
$ws.Range("E30").Value = '  -7.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0891'
$ws.Range("E31").Value = '  -8.44%  '
$ws.Range("E32").Value = '  -10.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '159.27'
$ws.Range("E33").Value = '  -7.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.02'
$ws.Range("E34").Value = '  -15.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.73'
$ws.Range("E35").Value = '  -7.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.08'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.123'
$ws.Range("E37").Value = '  -7.50%  '
$ws.Range("E40").Value = '  -12.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.56'
$ws.Range("E41").Value = '  -14.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0323'
$ws.Range("E42").Value = '  -10.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.816.50'
$ws.Range("E44").Value = '  +9.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.47'
$ws.Range("E45").Value = '  -11.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.05'
$ws.Range("E46").Value = '  -11.14%  '
$ws.Range("E47").Value = '  -14.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '77.41'
$ws.Range("E48").Value = '  -12.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.39'
$ws.Range("E49").Value = '  -5.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '60.54'
$ws.Range("E50").Value = '  -16.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.58'
$ws.Range("E51").Value = '  -9.44%  '

# Rows 38 and 39 swap content (ARBITRUM <-> RenderToken), with new price/volume values
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.48"
$ws.Range("E38").Value = "  -9.33%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.85"
$ws.Range("E39").Value = "  +3.80%  "
